$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.201709032058716
$ws.Range("B1").Value = 1.948946118354797
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.964933753013611
$ws.Range("E1").Value = 1.205845594406128
